$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Reln"
$ws.Range("C2").Value = "Lrp8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05950533333333333
$ws.Range("H2").Value = 0.178516
$ws.Range("I2").Value = 0.01803537323915772
$ws.Range("J2").Value = 0.01803537323915772
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.664391
$ws.Range("N2").Value = 4.993173000000001
$ws.Range("O2").Value = 0.3990511495040125
$ws.Range("P2").Value = 0.3990511495040125
$ws.Range("Q2").Value = 0.09904014125200002
$ws.Range("R2").Value = 0.8913612712680001
$ws.Range("S2").Value = 0.007197036422819795
$ws.Range("T2").Value = 0.007197036422819795

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Reln"
$ws.Range("C3").Value = "Lrp8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.05950533333333333
$ws.Range("H3").Value = 0.178516
$ws.Range("I3").Value = 0.01803537323915772
$ws.Range("J3").Value = 0.01803537323915772
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.437958
$ws.Range("N3").Value = 1.313874
$ws.Range("O3").Value = 0.1050039584054939
$ws.Range("P3").Value = 0.1050039584054938
$ws.Range("Q3").Value = 0.026060836776
$ws.Range("R3").Value = 0.234547530984
$ws.Range("S3").Value = 0.001893785581432074
$ws.Range("T3").Value = 0.001893785581432074

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Reln"
$ws.Range("C4").Value = "Lrp8"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.05950533333333333
$ws.Range("H4").Value = 0.178516
$ws.Range("I4").Value = 0.01803537323915772
$ws.Range("J4").Value = 0.01803537323915772
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.439215333333333
$ws.Range("N4").Value = 4.317646
$ws.Range("O4").Value = 0.3450634695516061
$ws.Range("P4").Value = 0.3450634695516061
$ws.Range("Q4").Value = 0.08564098814844444
$ws.Range("R4").Value = 0.770768893336
$ws.Range("S4").Value = 0.006223348464561951
$ws.Range("T4").Value = 0.006223348464561951

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Reln"
$ws.Range("C5").Value = "Lrp8"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.05950533333333333
$ws.Range("H5").Value = 0.178516
$ws.Range("I5").Value = 0.01803537323915772
$ws.Range("J5").Value = 0.01803537323915772
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.629307
$ws.Range("N5").Value = 1.887921
$ws.Range("O5").Value = 0.1508814225388875
$ws.Range("P5").Value = 0.1508814225388875
$ws.Range("Q5").Value = 0.037447122804
$ws.Range("R5").Value = 0.337024105236
$ws.Range("S5").Value = 0.002721202770343901
$ws.Range("T5").Value = 0.002721202770343901

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Reln"
$ws.Range("C6").Value = "Lrp8"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.06084700000000001
$ws.Range("H6").Value = 0.182541
$ws.Range("I6").Value = 0.01844201677412159
$ws.Range("J6").Value = 0.01844201677412159
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.664391
$ws.Range("N6").Value = 4.993173000000001
$ws.Range("O6").Value = 0.3990511495040125
$ws.Range("P6").Value = 0.3990511495040125
$ws.Range("Q6").Value = 0.101273199177
$ws.Range("R6").Value = 0.9114587925930001
$ws.Range("S6").Value = 0.007359307992885503
$ws.Range("T6").Value = 0.007359307992885503

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Reln"
$ws.Range("C7").Value = "Lrp8"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.06084700000000001
$ws.Range("H7").Value = 0.182541
$ws.Range("I7").Value = 0.01844201677412159
$ws.Range("J7").Value = 0.01844201677412159
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.437958
$ws.Range("N7").Value = 1.313874
$ws.Range("O7").Value = 0.1050039584054939
$ws.Range("P7").Value = 0.1050039584054938
$ws.Range("Q7").Value = 0.026648430426
$ws.Range("R7").Value = 0.239835873834
$ws.Range("S7").Value = 0.001936484762263284
$ws.Range("T7").Value = 0.001936484762263283

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Reln"
$ws.Range("C8").Value = "Lrp8"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.06084700000000001
$ws.Range("H8").Value = 0.182541
$ws.Range("I8").Value = 0.01844201677412159
$ws.Range("J8").Value = 0.01844201677412159
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.439215333333333
$ws.Range("N8").Value = 4.317646
$ws.Range("O8").Value = 0.3450634695516061
$ws.Range("P8").Value = 0.3450634695516061
$ws.Range("Q8").Value = 0.08757193538733334
$ws.Range("R8").Value = 0.788147418486
$ws.Range("S8").Value = 0.006363666293607314
$ws.Range("T8").Value = 0.006363666293607314

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Reln"
$ws.Range("C9").Value = "Lrp8"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.06084700000000001
$ws.Range("H9").Value = 0.182541
$ws.Range("I9").Value = 0.01844201677412159
$ws.Range("J9").Value = 0.01844201677412159
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.629307
$ws.Range("N9").Value = 1.887921
$ws.Range("O9").Value = 0.1508814225388875
$ws.Range("P9").Value = 0.1508814225388875
$ws.Range("Q9").Value = 0.038291443029
$ws.Range("R9").Value = 0.344622987261
$ws.Range("S9").Value = 0.002782557725365491
$ws.Range("T9").Value = 0.002782557725365491

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Reln"
$ws.Range("C10").Value = "Lrp8"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.179015666666666
$ws.Range("H10").Value = 9.537047
$ws.Range("I10").Value = 0.9635226099867207
$ws.Range("J10").Value = 0.9635226099867207
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.664391
$ws.Range("N10").Value = 4.993173000000001
$ws.Range("O10").Value = 0.3990511495040125
$ws.Range("P10").Value = 0.3990511495040125
$ws.Range("Q10").Value = 5.291125064459
$ws.Range("R10").Value = 47.620125580131
$ws.Range("S10").Value = 0.3844948050883072
$ws.Range("T10").Value = 0.3844948050883072

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Reln"
$ws.Range("C11").Value = "Lrp8"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.179015666666666
$ws.Range("H11").Value = 9.537047
$ws.Range("I11").Value = 0.9635226099867207
$ws.Range("J11").Value = 0.9635226099867207
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.437958
$ws.Range("N11").Value = 1.313874
$ws.Range("O11").Value = 0.1050039584054939
$ws.Range("P11").Value = 0.1050039584054938
$ws.Range("Q11").Value = 1.392275343342
$ws.Range("R11").Value = 12.530478090078
$ws.Range("S11").Value = 0.1011736880617985
$ws.Range("T11").Value = 0.1011736880617985

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Reln"
$ws.Range("C12").Value = "Lrp8"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.179015666666666
$ws.Range("H12").Value = 9.537047
$ws.Range("I12").Value = 0.9635226099867207
$ws.Range("J12").Value = 0.9635226099867207
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.439215333333333
$ws.Range("N12").Value = 4.317646
$ws.Range("O12").Value = 0.3450634695516061
$ws.Range("P12").Value = 0.3450634695516061
$ws.Range("Q12").Value = 4.575288092373555
$ws.Range("R12").Value = 41.177592831362
$ws.Range("S12").Value = 0.3324764547934368
$ws.Range("T12").Value = 0.3324764547934368

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Reln"
$ws.Range("C13").Value = "Lrp8"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.179015666666666
$ws.Range("H13").Value = 9.537047
$ws.Range("I13").Value = 0.9635226099867207
$ws.Range("J13").Value = 0.9635226099867207
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.629307
$ws.Range("N13").Value = 1.887921
$ws.Range("O13").Value = 0.1508814225388875
$ws.Range("P13").Value = 0.1508814225388875
$ws.Range("Q13").Value = 2.000576812143
$ws.Range("R13").Value = 18.005191309287
$ws.Range("S13").Value = 0.3324764547934368
$ws.Range("T13").Value = 0.3324764547934368

Write-Output "edit complete"